$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(5,1).Value = 112241869
$ws.Cells.Item(5,2).Value = 89539
$ws.Cells.Item(5,17).Value = 553972
$ws.Cells.Item(5,18).Value = 7008047
$ws.Cells.Item(6,1).Value = 112241883
$ws.Cells.Item(6,2).Value = 78699
$ws.Cells.Item(6,4).Value = "NT"
$ws.Cells.Item(6,5).Value = 6458
$ws.Cells.Item(6,6).Value = "Lunglav"
$ws.Cells.Item(6,7).Value = "Lobaria pulmonaria"
$ws.Cells.Item(6,8).Value = "(L.) Hoffm."
$ws.Cells.Item(6,17).Value = 553961
$ws.Cells.Item(6,18).Value = 7008062
$ws.Cells.Item(7,1).Value = 112241882
$ws.Cells.Item(7,2).Value = 78699
$ws.Cells.Item(7,5).Value = 6458
$ws.Cells.Item(7,6).Value = "Lunglav"
$ws.Cells.Item(7,7).Value = "Lobaria pulmonaria"
$ws.Cells.Item(7,8).Value = "(L.) Hoffm."
$ws.Cells.Item(7,17).Value = 553964
$ws.Cells.Item(7,18).Value = 7008000
$ws.Cells.Item(8,1).Value = 112241881
$ws.Cells.Item(8,2).Value = 78699
$ws.Cells.Item(8,5).Value = 6458
$ws.Cells.Item(8,6).Value = "Lunglav"
$ws.Cells.Item(8,7).Value = "Lobaria pulmonaria"
$ws.Cells.Item(8,8).Value = "(L.) Hoffm."
$ws.Cells.Item(8,17).Value = 554185
$ws.Cells.Item(8,18).Value = 7007926
$ws.Cells.Item(9,2).Value = 78699
$ws.Cells.Item(10,1).Value = 112241876
$ws.Cells.Item(10,2).Value = 56446
$ws.Cells.Item(10,5).Value = 100049
$ws.Cells.Item(10,6).Value = "Spillkråka"
$ws.Cells.Item(10,7).Value = "Dryocopus martius"
$ws.Cells.Item(10,8).Value = "(Linnaeus, 1758)"
$ws.Cells.Item(10,9).Value = "1"
$ws.Cells.Item(10,13).Value = "födosökande"
$ws.Cells.Item(10,14).Value = "observerad"
$ws.Cells.Item(10,17).Value = 554066
$ws.Cells.Item(10,18).Value = 7008014
$ws.Cells.Item(11,1).Value = 112241880
$ws.Cells.Item(11,2).Value = 78699
$ws.Cells.Item(11,5).Value = 6458
$ws.Cells.Item(11,6).Value = "Lunglav"
$ws.Cells.Item(11,7).Value = "Lobaria pulmonaria"
$ws.Cells.Item(11,8).Value = "(L.) Hoffm."
$ws.Cells.Item(11,17).Value = 554164
$ws.Cells.Item(11,18).Value = 7007925
$ws.Cells.Item(12,1).Value = 112241867
$ws.Cells.Item(12,2).Value = 89539
$ws.Cells.Item(12,5).Value = 1202
$ws.Cells.Item(12,6).Value = "Ullticka"
$ws.Cells.Item(12,7).Value = "Phellinidium ferrugineofuscum"
$ws.Cells.Item(12,8).Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Cells.Item(12,9).Value = ""
$ws.Cells.Item(12,13).Value = ""
$ws.Cells.Item(12,14).Value = ""
$ws.Cells.Item(12,17).Value = 554177
$ws.Cells.Item(12,18).Value = 7007859
$ws.Cells.Item(13,1).Value = 112241871
$ws.Cells.Item(13,2).Value = 89539
$ws.Cells.Item(13,5).Value = 1202
$ws.Cells.Item(13,6).Value = "Ullticka"
$ws.Cells.Item(13,7).Value = "Phellinidium ferrugineofuscum"
$ws.Cells.Item(13,8).Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Cells.Item(13,17).Value = 554086
$ws.Cells.Item(13,18).Value = 7008007
$ws.Cells.Item(14,2).Value = 78699
$ws.Cells.Item(15,1).Value = 112241868
$ws.Cells.Item(15,2).Value = 89539
$ws.Cells.Item(15,5).Value = 1202
$ws.Cells.Item(15,6).Value = "Ullticka"
$ws.Cells.Item(15,7).Value = "Phellinidium ferrugineofuscum"
$ws.Cells.Item(15,8).Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Cells.Item(15,17).Value = 554007
$ws.Cells.Item(15,18).Value = 7007988
$ws.Cells.Item(16,1).Value = 112241870
$ws.Cells.Item(16,2).Value = 89539
$ws.Cells.Item(16,5).Value = 1202
$ws.Cells.Item(16,6).Value = "Ullticka"
$ws.Cells.Item(16,7).Value = "Phellinidium ferrugineofuscum"
$ws.Cells.Item(16,8).Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Cells.Item(16,17).Value = 554078
$ws.Cells.Item(16,18).Value = 7008009
$ws.Cells.Item(17,1).Value = 112241877
$ws.Cells.Item(17,2).Value = 89979
$ws.Cells.Item(17,4).Value = "VU"
$ws.Cells.Item(17,5).Value = 1209
$ws.Cells.Item(17,6).Value = "Rynkskinn"
$ws.Cells.Item(17,7).Value = "Phlebia centrifuga"
$ws.Cells.Item(17,8).Value = "P.Karst."
$ws.Cells.Item(17,17).Value = 554081
$ws.Cells.Item(17,18).Value = 7007966
$ws.Cells.Item(18,1).Value = 112304906
$ws.Cells.Item(18,2).Value = 89485
$ws.Cells.Item(18,3).Value = "Ovaliderad"
$ws.Cells.Item(18,4).Value = "NT"
$ws.Cells.Item(18,5).Value = 112
$ws.Cells.Item(18,6).Value = "Stjärntagging"
$ws.Cells.Item(18,7).Value = "Asterodon ferruginosus"
$ws.Cells.Item(18,8).Value = "Pat."
$ws.Cells.Item(18,16).Value = "sotberget, Jmt"
$ws.Cells.Item(18,17).Value = 553996
$ws.Cells.Item(18,18).Value = 7008025
$ws.Cells.Item(18,19).Value = 10
$ws.Cells.Item(18,20).Value = "Jämtland"
$ws.Cells.Item(18,21).Value = "Ragunda"
$ws.Cells.Item(18,22).Value = "Jämtland"
$ws.Cells.Item(18,23).Value = "Ragunda"
$ws.Cells.Item(18,25).Value = "2023-09-21"
$ws.Cells.Item(18,27).Value = "2023-09-21"
$ws.Cells.Item(18,30).Value = $false
$ws.Cells.Item(18,31).Value = $false
$ws.Cells.Item(18,33).Value = $false
$ws.Cells.Item(18,49).Value = "Benny Öwre"
$ws.Cells.Item(18,50).Value = "Benny Öwre"
$ws.Cells.Item(19,1).Value = 112304898
$ws.Cells.Item(19,2).Value = 90226
$ws.Cells.Item(19,3).Value = "Ovaliderad"
$ws.Cells.Item(19,4).Value = "VU"
$ws.Cells.Item(19,5).Value = 67
$ws.Cells.Item(19,6).Value = "Sprickporing"
$ws.Cells.Item(19,7).Value = "Diplomitoporus crustulinus"
$ws.Cells.Item(19,8).Value = "(Bres.) Domański"
$ws.Cells.Item(19,16).Value = "Sotberget, Jmt"
$ws.Cells.Item(19,17).Value = 553981
$ws.Cells.Item(19,18).Value = 7008031
$ws.Cells.Item(19,19).Value = 10
$ws.Cells.Item(19,20).Value = "Jämtland"
$ws.Cells.Item(19,21).Value = "Ragunda"
$ws.Cells.Item(19,22).Value = "Jämtland"
$ws.Cells.Item(19,23).Value = "Ragunda"
$ws.Cells.Item(19,25).Value = "2023-09-21"
$ws.Cells.Item(19,27).Value = "2023-09-21"
$ws.Cells.Item(19,30).Value = $false
$ws.Cells.Item(19,31).Value = $false
$ws.Cells.Item(19,33).Value = $false
$ws.Cells.Item(19,49).Value = "Benny Öwre"
$ws.Cells.Item(19,50).Value = "Benny Öwre"
